$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rocket dimension measurements + Center-of-Pressure (CP) calculations added
# to the right of the existing mass table (columns G:J).
# ---------------------------------------------------------------------------

# ---- Normal Force Coefficient section -------------------------------------
$ws.Range("G1").Value = "Normal Force Coefficient"

$ws.Range("G2").Value = "NoseCone"
$ws.Range("H2").Value = "CNn"
$ws.Range("I2").Value = 2

$ws.Range("G4").Value = "finset (4)"
$ws.Range("H4").Value = "f"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = "m"

$ws.Range("H5").Value = "R"
$ws.Range("I5").Formula = "=3.5/2/100"
$ws.Range("J5").Value = "m"

$ws.Range("H6").Value = "s"
$ws.Range("I6").Formula = "=5/100"
$ws.Range("J6").Value = "m"

$ws.Range("H7").Value = "N"
$ws.Range("I7").Value = 4

$ws.Range("H8").Value = "d"
$ws.Range("I8").Formula = "=3.5/100"
$ws.Range("J8").Value = "m"

$ws.Range("H9").Value = "xr"
$ws.Range("I9").Formula = "=4/100"
$ws.Range("J9").Value = "m"

$ws.Range("H10").Value = "ct"
$ws.Range("I10").Formula = "=3/100"
$ws.Range("J10").Value = "m"

$ws.Range("H11").Value = "cr"
$ws.Range("I11").Formula = "=5.5/100"
$ws.Range("J11").Value = "m"

$ws.Range("H12").Value = "theta"
$ws.Range("I12").Formula = "=ATAN((1/I6)*(I9+0.5*(I10-I11)))"
$ws.Range("J12").Value = "rad"

$ws.Range("H13").Value = "l"
$ws.Range("I13").Formula = "=I6/COS(I12)"

$ws.Range("H15").Value = "CNf"
$ws.Range("I15").Formula = "=(1+(I4*I5)/(I6+I5))*((4*I7*(I6/I8)^2)/(1+SQRT(1+((2*I13)/(I11+I10))^2)))"

# ---- Boattail section -------------------------------------------------------
$ws.Range("G17").Value = "Boattail"
$ws.Range("H17").Value = "S1"
$ws.Range("I17").Formula = "=1/4*PI()*I8^2"
$ws.Range("J17").Value = "m^2"

$ws.Range("H18").Value = "S2"
$ws.Range("I18").Formula = "=1/4*PI()*0.03^2"
$ws.Range("J18").Value = "m^2"

$ws.Range("H20").Value = "CNcb"
$ws.Range("I20").Formula = "=(8/(PI()*I8^2))*(I18-I17)"

# ---- Center of Pressure Location section -----------------------------------
$ws.Range("G22").Value = "Center of Pressure Location"

$ws.Range("G23").Value = "Nosecone (Ogive)"
$ws.Range("H23").Value = "Xn"
$ws.Range("I23").Formula = "=0.466*0.06"

$ws.Range("G25").Value = "Finset"
$ws.Range("H25").Value = "XB"
$ws.Range("I25").Value = 0.68

$ws.Range("H26").Value = "XF"
$ws.Range("I26").Formula = "=I25+I9/3*(I11+2*I10)/(I11+I10)+1/6*(I11+I10-(I11*I10)/(I11+I10))"

$ws.Range("G28").Value = "Conical Transitional"
$ws.Range("H28").Value = "LT"
$ws.Range("I28").Formula = "=1.2/100"

$ws.Range("H29").Value = "XP"
$ws.Range("I29").Value = 0.735

$ws.Range("H30").Value = "XT"
$ws.Range("I30").Formula = "=I29+I28/3*(1+(1-I8/0.03)/(1-(I8/0.03)^2))"

$ws.Range("G32").Value = "Total CP"
$ws.Range("H32").Value = "X"
$ws.Range("I32").Formula = "=(I2*I23+I15*I26+I20*I30)/(I2+I15+I20)"
$ws.Range("J32").Value = "m"

# ---------------------------------------------------------------------------
# Formatting: bold header/result cells.
# ---------------------------------------------------------------------------

# Title matches the existing bold "Component"/"Mass (g)" header style.
$ws.Range("G1").Font.Bold = $true

# Bold "result" cells get their own (new) bold font entry.
$rBold = $ws.Range("H2,I2,H15,I15,H20,I20,G22,H23,I23,H26,I26,H30,I30,H32,I32")
$rBold.Font.Bold = $true
$rBold.Font.Name = "Aptos Narrow"

# H29 ("XP") carries its own distinct (non-bold) font entry.
$ws.Range("H29").Font.Name = "Aptos Narrow"

# ---------------------------------------------------------------------------
# Column widths for the new columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 25.666666666666668
$ws.Columns.Item(8).ColumnWidth = 4.666666666666667
$ws.Columns.Item(9).ColumnWidth = 11.041666666666666

# ---------------------------------------------------------------------------
# Sheet view: scrolled / selection state.
# ---------------------------------------------------------------------------
$ws.Range("I35").Select()
